# Added HitRate into lowcode.RTP
# This re-sorts / re-assigns the data rows (A2:F25) of the active sheet
# to reflect the new "HitRate" ordering used by lowcode.RTP.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for rows 2-25, columns A-F (symbol, reel1..reel5)
$data = @(
    @(901,  16, 15, 45, 60, 60),
    @(301,  6,  45, 30, 60, 45),
    @(1201, 2,  10, 10, 10, 10),
    @(1202, 2,  10, 10, 10, 10),
    @(1203, 3,  15, 15, 15, 15),
    @(801,  3,  67, 65, 52, 45),
    @(902,  1,  0,  0,  0,  0),
    @(601,  9,  60, 67, 60, 42),
    @(201,  9,  30, 15, 45, 30),
    @(101,  9,  30, 15, 60, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(501,  9,  52, 30, 75, 45),
    @(401,  9,  48, 67, 75, 45),
    @(701,  3,  90, 45, 97, 15),
    @(502,  0,  4,  0,  0,  0),
    @(802,  0,  4,  5,  4,  0),
    @(2,    0,  2,  2,  2,  2),
    @(3,    0,  3,  3,  3,  3),
    @(1101, 0,  15, 30, 30, 0),
    @(1,    0,  2,  2,  2,  2),
    @(402,  0,  0,  4,  0,  0),
    @(602,  0,  0,  4,  0,  9),
    @(702,  0,  0,  0,  4,  0),
    @(1002, 0,  0,  0,  0,  9)
)

$startRow = 2
$cols = @("A", "B", "C", "D", "E", "F")
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + "$row").Value = $values[$j]
    }
}
